$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Cxcl1"
$ws.Range("C2").Value = "Cxcr1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 81.50654200000001
$ws.Range("H2").Value = 244.519626
$ws.Range("I2").Value = 0.2392262902761414
$ws.Range("J2").Value = 0.2392262902761414
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 2.161116666666667
$ws.Range("N2").Value = 6.48335
$ws.Range("O2").Value = 0.9584259542460944
$ws.Range("P2").Value = 0.9584259542460944
$ws.Range("Q2").Value = 176.1451463585667
$ws.Range("R2").Value = 1585.3063172271
$ws.Range("S2").Value = 0.229280685538664
$ws.Range("T2").Value = 0.229280685538664

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Cxcl1"
$ws.Range("C3").Value = "Cxcr1"
$ws.Range("D3").Value = "M2"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 81.50654200000001
$ws.Range("H3").Value = 244.519626
$ws.Range("I3").Value = 0.2392262902761414
$ws.Range("J3").Value = 0.2392262902761414
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.09374366666666666
$ws.Range("N3").Value = 0.281231
$ws.Range("O3").Value = 0.04157404575390553
$ws.Range("P3").Value = 0.04157404575390552
$ws.Range("Q3").Value = 7.640722104400667
$ws.Range("R3").Value = 68.766498939606
$ws.Range("S3").Value = 0.009945604737477388
$ws.Range("T3").Value = 0.009945604737477387

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Cxcl1"
$ws.Range("C4").Value = "Cxcr1"
$ws.Range("D4").Value = "ECs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 195.050573
$ws.Range("H4").Value = 585.1517190000001
$ws.Range("I4").Value = 0.5724844147482752
$ws.Range("J4").Value = 0.5724844147482753
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 2.161116666666667
$ws.Range("N4").Value = 6.48335
$ws.Range("O4").Value = 0.9584259542460944
$ws.Range("P4").Value = 0.9584259542460944
$ws.Range("Q4").Value = 421.5270441531833
$ws.Range("R4").Value = 3793.74339737865
$ws.Range("S4").Value = 0.5486839214961325
$ws.Range("T4").Value = 0.5486839214961327

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Cxcl1"
$ws.Range("C5").Value = "Cxcr1"
$ws.Range("D5").Value = "M2"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 195.050573
$ws.Range("H5").Value = 585.1517190000001
$ws.Range("I5").Value = 0.5724844147482752
$ws.Range("J5").Value = 0.5724844147482753
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.09374366666666666
$ws.Range("N5").Value = 0.281231
$ws.Range("O5").Value = 0.04157404575390553
$ws.Range("P5").Value = 0.04157404575390552
$ws.Range("Q5").Value = 18.28475589845433
$ws.Range("R5").Value = 164.562803086089
$ws.Range("S5").Value = 0.02380049325214262
$ws.Range("T5").Value = 0.02380049325214262

# Row 6 (new)
$ws.Range("A6").Value = "M2"
$ws.Range("B6").Value = "Cxcl1"
$ws.Range("C6").Value = "Cxcr1"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 19.93351866666666
$ws.Range("H6").Value = 59.80055599999999
$ws.Range("I6").Value = 0.05850599971198487
$ws.Range("J6").Value = 0.05850599971198487
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 2.161116666666667
$ws.Range("N6").Value = 6.48335
$ws.Range("O6").Value = 0.9584259542460944
$ws.Range("P6").Value = 0.9584259542460944
$ws.Range("Q6").Value = 43.07865941584443
$ws.Range("R6").Value = 387.7079347425999
$ws.Range("S6").Value = 0.05607366860308083
$ws.Range("T6").Value = 0.05607366860308083

# Row 7 (new)
$ws.Range("A7").Value = "M2"
$ws.Range("B7").Value = "Cxcl1"
$ws.Range("C7").Value = "Cxcr1"
$ws.Range("D7").Value = "M2"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 19.93351866666666
$ws.Range("H7").Value = 59.80055599999999
$ws.Range("I7").Value = 0.05850599971198487
$ws.Range("J7").Value = 0.05850599971198487
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.09374366666666666
$ws.Range("N7").Value = 0.281231
$ws.Range("O7").Value = 0.04157404575390553
$ws.Range("P7").Value = 0.04157404575390552
$ws.Range("Q7").Value = 1.868641129381777
$ws.Range("R7").Value = 16.81777016443599
$ws.Range("S7").Value = 0.002432331108904042
$ws.Range("T7").Value = 0.002432331108904042

# Row 8 (new)
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Cxcl1"
$ws.Range("C8").Value = "Cxcr1"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 44.21833233333334
$ws.Range("H8").Value = 132.654997
$ws.Range("I8").Value = 0.1297832952635985
$ws.Range("J8").Value = 0.1297832952635985
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 2.161116666666667
$ws.Range("N8").Value = 6.48335
$ws.Range("O8").Value = 0.9584259542460944
$ws.Range("P8").Value = 0.9584259542460944
$ws.Range("Q8").Value = 95.56097497777222
$ws.Range("R8").Value = 860.04877479995
$ws.Range("S8").Value = 0.124387678608217
$ws.Range("T8").Value = 0.124387678608217

# Row 9 (new)
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Cxcl1"
$ws.Range("C9").Value = "Cxcr1"
$ws.Range("D9").Value = "M2"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 44.21833233333334
$ws.Range("H9").Value = 132.654997
$ws.Range("I9").Value = 0.1297832952635985
$ws.Range("J9").Value = 0.1297832952635985
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.09374366666666666
$ws.Range("N9").Value = 0.281231
$ws.Range("O9").Value = 0.04157404575390553
$ws.Range("P9").Value = 0.04157404575390552
$ws.Range("Q9").Value = 4.145188606811889
$ws.Range("R9").Value = 37.306697461307
$ws.Range("S9").Value = 0.005395616655381473
$ws.Range("T9").Value = 0.005395616655381472
